$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a "Population/Dwellings" block at the bottom:
#   Row 27: Population Total
#   Row 28: Attached Dwellings
#   Row 29: Houses
# We need to insert 4 new demographic rows (Net Interstate Migration,
# Net Overseas Migration, Natural Increase, Total Population Increase)
# right after "Population Total" (row 27), pushing the Attached
# Dwellings / Houses rows down to 32/33.

$ws.Rows("28:31").Insert()

# Excel's row-insert copies the formatting of the row above (row 27,
# which uses the "odd" banding style) into all four new rows. Fix the
# banding so it alternates correctly (row28/30 = "even" style like the
# Attached Dwellings row now sitting at row32, row29/31 = "odd" style).
$ws.Range("A32:L32").Copy()
$ws.Range("A28:L28").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:L30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the row labels (column B) top-to-bottom first ...
$ws.Range("B28").Value2 = "Net Interstate Migration"
$ws.Range("B29").Value2 = "Net Overseas Migration"
$ws.Range("B30").Value2 = "Natural Increase"
$ws.Range("B31").Value2 = "Total Population Increase"

# ... then the series codes (column C), entered in the order
# 28, 29, 31, 30 to reproduce the authoring sequence.
$ws.Range("C28").Value2 = "NIMTOT"
$ws.Range("C29").Value2 = "NOMTOT"
$ws.Range("C31").Value2 = "POPINC"
$ws.Range("C30").Value2 = "NATTOT"

# Fill the rebasing coefficient data (all 1s) for the new rows.
$cols = @("D","E","F","G","H","I","J","K","L")
foreach ($r in 28..31) {
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value2 = 1
    }
}

# Match the saved selection shown in the workbook after editing.
$ws.Range("C30").Select()
